# ============================================================================
# Restructure PlayerPerformance_4581.xlsx:
#   - the single existing sheet ("ODI Batting") becomes "Player Info"
#     (new ID/NAME/BATTING_HAND/BOWL_STYLE content)
#   - a new "ODI Batting" sheet is inserted holding the original batting
#     log, with MATCH_CARD_LINK (a full URL) replaced by MATCH_CODE (just
#     the numeric code)
#   - a new "ODI Batting Extra" sheet is appended with additional
#     per-match batting stats
# ============================================================================

$wb = $excel.ActiveWorkbook
$excel.ScreenUpdating = $false

$original = $wb.Worksheets.Item(1)   # currently "ODI Batting" -- holds the batting log

# --- Step 1: clone the original batting log (values + formats) BEFORE we --
#     touch the original sheet at all, onto a throwaway-named sheet (the
#     original still owns the name "ODI Batting" at this point, so the real
#     target sheet has to be named after the rename in step 2).
$battingSheet = $wb.Worksheets.Add($null, $original)
$battingSheet.Name = "ODI Batting (new)"

$extraSheet = $wb.Worksheets.Add($null, $battingSheet)
$extraSheet.Name = "ODI Batting Extra"

$original.UsedRange.Copy($battingSheet.Range("A1"))
$excel.CutCopyMode = $false

# Grab a copy of the standard header look (bold font, thin border, centred)
# from the cloned sheet before the original's own copy of it gets wiped.
$headerStyleSource = $battingSheet.Range("A1")

# --- Step 2: free up the "ODI Batting" name by turning the original sheet -
#     into "Player Info", then claim the name for the clone.
$original.Cells.Clear()
$original.Name = "Player Info"
$battingSheet.Name = "ODI Batting"

# --- Step 3: rewrite the MATCH_CARD_LINK column on the clone into MATCH_CODE
$lastRow = $battingSheet.Cells.Item($battingSheet.Rows.Count, 1).End(-4162).Row

$battingSheet.Range("D1").Value = "MATCH_CODE"

$codeRange = $battingSheet.Range("D2:D" + $lastRow)
$codeRange.NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $url = $battingSheet.Cells.Item($r, 4).Value2
    $code = $url.Substring($url.IndexOf("MatchCode=") + 10)
    $battingSheet.Cells.Item($r, 4).Value = $code
}

# --- Step 4: populate "Player Info" with the player's bio row -------------
# Reuse the existing bold/bordered/centred header style (captured off the
# batting sheet's own header cell) instead of re-deriving it by hand, so we
# don't fork a second, visually-identical style record.
$headerStyleSource.Copy()
$original.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le 4; $c++) {
    $original.Cells.Item(1, $c).Value = $piHeaders[$c - 1]
}

$original.Range("A2").NumberFormat = "@"
$piRow = @("4581", "Weerahandige Inol Avishka Fernando", "Right Handed", "Right Arm Medium Fast")
for ($c = 1; $c -le 4; $c++) {
    $original.Cells.Item(2, $c).Value = $piRow[$c - 1]
}

# --- Step 5: populate "ODI Batting Extra" with the per-match extras -------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le 6; $c++) {
    $extraSheet.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
$extraHeaderRange = $extraSheet.Range("A1:F1")
$extraHeaderRange.Font.Bold = $true
$extraHeaderRange.HorizontalAlignment = -4108
$extraHeaderRange.VerticalAlignment = -4160
$extraHeaderRange.Borders.LineStyle = 1

# MATCH_CODE (A) and the numeric-looking text columns (C, D, E) must stay
# text, like the source data -- format them before writing.
$extraSheet.Range("A2:A21").NumberFormat = "@"
$extraSheet.Range("C2:E21").NumberFormat = "@"

$extraRows = @(
    @("4350", 3,    "4",  "0",  "7.58%",  "NO"),
    @("4356", $null, $null, $null, $null, "NO"),
    @("4357", $null, $null, $null, $null, "NO"),
    @("4358", $null, $null, $null, $null, "NO"),
    @("4375", 3,    "0",  "0",  $null,    "NO"),
    @("4376", 2,    "1",  "0",  "1.35%",  "NO"),
    @("4413", 1,    "5",  "1",  "17.24%", "NO"),
    @("4414", $null, $null, $null, $null, "NO"),
    @("4417", 1,    "5",  "0",  "9.45%",  "NO"),
    @("4470", 3,    "0",  "0",  "0.83%",  "NO"),
    @("4471", $null, $null, $null, $null, "NO"),
    @("4480", 1,    "2",  "1",  "12.60%", "NO"),
    @("4482", 1,    "4",  "1",  "18.18%", "NO"),
    @("4485", 1,    "4",  "1",  "33.48%", "YES"),
    @("4487", 1,    "10", "2",  "39.33%", "YES"),
    @("4488", 1,    "1",  "0",  "4.06%",  "NO"),
    @("4491", 1,    "2",  "0",  "4.93%",  "NO"),
    @("4687", $null, $null, $null, $null, "NO"),
    @("4689", 1,    "4",  "0",  "9.30%",  "NO"),
    @("4691", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extraSheet.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) { $extraSheet.Cells.Item($r, 2).Value = $row[1] }
    if ($null -ne $row[2]) { $extraSheet.Cells.Item($r, 3).Value = $row[2] }
    if ($null -ne $row[3]) { $extraSheet.Cells.Item($r, 4).Value = $row[3] }
    if ($null -ne $row[4]) { $extraSheet.Cells.Item($r, 5).Value = $row[4] }
    $extraSheet.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Step 6: leave the workbook selection on the first sheet, as before ---
$original.Activate()
[void]$original.Range("A1").Select()
$excel.ScreenUpdating = $true

Write-Output "Parts 1-5 done"
